$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.363.40'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '3.026.28'
$ws.Range('E3').Value = '  -4.64%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''567.51'
$ws.Range('E5').Value = '  -3.87%  '
$ws.Range('D6').Value = '''129.69'
$ws.Range('E6').Value = '  -5.10%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.022.00'
$ws.Range('E8').Value = '  -4.69%  '
$ws.Range('D9').Value = '''0.500'
$ws.Range('E9').Value = '  -2.35%  '
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('D11').Value = '''5.30'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  -4.85%  '
$ws.Range('E13').Value = '  -5.22%  '
$ws.Range('D14').Value = '''33.28'
$ws.Range('E14').Value = '  -4.52%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '3.525.56'
$ws.Range('E16').Value = '  -4.60%  '
$ws.Range('D17').Value = '61.461.24'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '3.014.60'
$ws.Range('E18').Value = '  -4.96%  '
$ws.Range('D19').Value = '''6.24'
$ws.Range('E19').Value = '  -5.45%  '
$ws.Range('D20').Value = '''441.24'
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('D21').Value = '''13.25'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').Value = '''0.667'
$ws.Range('E22').Value = '  -6.26%  '
$ws.Range('E23').Value = '  -5.93%  '
$ws.Range('D24').Value = '''12.99'
$ws.Range('E24').Value = '  -3.22%  '
$ws.Range('D25').Value = '''79.47'
$ws.Range('E25').Value = '  -4.65%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -6.52%  '
$ws.Range('D29').Value = '''7.29'
$ws.Range('E29').Value = '  -6.25%  '
$ws.Range('E30').Value = '  -6.54%  '
$ws.Range('D31').Value = '''6.19'
$ws.Range('E31').Value = '  -9.45%  '
$ws.Range('D32').Value = '''25.64'
$ws.Range('E32').Value = '  -6.36%  '
$ws.Range('D33').Value = '''0.0945'
$ws.Range('E33').Value = '  -8.70%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('D35').Value = '''0.965'
$ws.Range('E35').Value = '  -7.16%  '
$ws.Range('D36').Value = '''5.64'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').Value = '''50.43'
$ws.Range('D38').Value = '0.0₃0678'
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('D39').Value = '''0.0363'
$ws.Range('E39').Value = '  -6.29%  '
$ws.Range('E40').Value = '  -3.63%  '
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').Value = '''382.37'
$ws.Range('E42').Value = '  -4.94%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.690.44'
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.48'
$ws.Range('E44').Value = '  -7.94%  '
$ws.Range('D46').Value = '''34.53'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('E47').Value = '  -6.05%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''120.16'
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '''1.99'
$ws.Range('E49').Value = '  -6.48%  '
$ws.Range('E50').Value = '  -3.68%  '
$ws.Range('D51').Value = '''23.51'
$ws.Range('E51').Value = '  -8.23%  '
